$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-05-17 Saturday"; new = "2025-05-18 Sunday"},
    @{old = "592×4="; new = "198×7="},
    @{old = "181×4="; new = "360×2="},
    @{old = "521×4="; new = "201×4="},
    @{old = "557×8="; new = "892×5="},
    @{old = "364×6="; new = "304×5="},
    @{old = "854×7="; new = "392×5="},
    @{old = "967×2="; new = "104×7="},
    @{old = "628×7="; new = "269×2="},
    @{old = "620×3="; new = "852×2="},
    @{old = "130×4="; new = "148×4="},
    @{old = "502×3="; new = "507×7="},
    @{old = "288×8="; new = "107×4="},
    @{old = "756×4="; new = "819×8="},
    @{old = "527×3="; new = "926×2="},
    @{old = "965×3="; new = "953×4="},
    @{old = "906×3="; new = "126×4="},
    @{old = "130×7="; new = "120×2="},
    @{old = "962×5="; new = "245×5="},
    @{old = "174×5="; new = "856×9="},
    @{old = "471×4="; new = "550×4="},
    @{old = "865×2="; new = "516×8="},
    @{old = "979×9="; new = "901×2="},
    @{old = "161×2="; new = "672×2="},
    @{old = "569×3="; new = "280×5="},
    @{old = "473×8="; new = "389×9="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
